# Apply odds updates scraped for Jogos_da_Semana_FlashScore_2025-04-15.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AB2").Value = 29
$ws.Range("AC2").Value = 151
$ws.Range("AF2").Value = 17
$ws.Range("AG2").Value = 17
$ws.Range("AH2").Value = 51
$ws.Range("AI2").Value = 51
$ws.Range("AJ2").Value = 81
$ws.Range("G2").Value = 2.3
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 1.2
$ws.Range("K2").Value = 4.33
$ws.Range("L2").Value = 1.91
$ws.Range("M2").Value = 1.8
$ws.Range("P2").Value = 1.93
$ws.Range("Q2").Value = 1.88
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 1.36
$ws.Range("T2").Value = 4.5
$ws.Range("U2").Value = 8.5
$ws.Range("W2").Value = 23
$ws.Range("X2").Value = 34
$ws.Range("Z2").Value = 4.33
# Row 4
$ws.Range("J4").Value = 1.08
$ws.Range("K4").Value = 8
# Row 5
$ws.Range("AB5").Value = 12
$ws.Range("AG5").Value = 9
$ws.Range("AJ5").Value = 21
$ws.Range("G5").Value = 3
$ws.Range("I5").Value = 2.15
$ws.Range("L5").Value = 1.17
$ws.Range("M5").Value = 5
$ws.Range("N5").Value = 1.53
$ws.Range("O5").Value = 2.4
$ws.Range("R5").Value = 1.5
$ws.Range("S5").Value = 2.5
$ws.Range("T5").Value = 13
$ws.Range("Z5").Value = 17
# Row 6
$ws.Range("AA6").Value = 7
$ws.Range("AC6").Value = 51
$ws.Range("AE6").Value = 7.5
$ws.Range("AF6").Value = 8
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 3.75
$ws.Range("L6").Value = 1.25
$ws.Range("M6").Value = 3.75
$ws.Range("N6").Value = 1.83
$ws.Range("O6").Value = 2.03
$ws.Range("P6").Value = 1.36
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 1.95
$ws.Range("Z6").Value = 12
# Row 7
$ws.Range("AA7").Value = 8
$ws.Range("AB7").Value = 13
$ws.Range("AE7").Value = 15
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 51
$ws.Range("G7").Value = 1.73
$ws.Range("H7").Value = 4.1
$ws.Range("L7").Value = 1.18
$ws.Range("M7").Value = 4.5
$ws.Range("N7").Value = 1.62
$ws.Range("O7").Value = 2.25
$ws.Range("P7").Value = 1.3
$ws.Range("Q7").Value = 3.4
$ws.Range("R7").Value = 1.62
$ws.Range("S7").Value = 2.2
$ws.Range("T7").Value = 9
$ws.Range("U7").Value = 9.5
$ws.Range("Z7").Value = 15
# Row 10
$ws.Range("J10").Value = 1.1
$ws.Range("K10").Value = 7
# Row 11
$ws.Range("J11").Value = 1.07
$ws.Range("K11").Value = 9
# Row 14
$ws.Range("J14").Value = 1.06
$ws.Range("K14").Value = 10
$ws.Range("L14").Value = 1.3
$ws.Range("M14").Value = 3.4
$ws.Range("N14").Value = 2
$ws.Range("O14").Value = 1.8
# Row 15
$ws.Range("AB15").Value = 15
$ws.Range("AC15").Value = 51
$ws.Range("AD15").Value = 251
# Row 16
$ws.Range("AA16").Value = 12
$ws.Range("AB16").Value = 30
$ws.Range("AC16").Value = 150
$ws.Range("AE16").Value = 6.4
$ws.Range("AF16").Value = 4.8
$ws.Range("AG16").Value = 9
$ws.Range("AH16").Value = 5.1
$ws.Range("AI16").Value = 9.5
$ws.Range("AJ16").Value = 32
$ws.Range("G16").Value = 17.5
$ws.Range("H16").Value = 6.1
$ws.Range("I16").Value = 1.13
$ws.Range("N16").Value = 1.45
$ws.Range("O16").Value = 2.55
$ws.Range("R16").Value = 2.5
$ws.Range("S16").Value = 1.47
$ws.Range("T16").Value = 40
$ws.Range("U16").Value = 150
$ws.Range("V16").Value = 50
$ws.Range("X16").Value = 250
$ws.Range("Y16").Value = 150
$ws.Range("Z16").Value = 14.5
# Row 17
$ws.Range("AA17").Value = 5.8
$ws.Range("AC17").Value = 40
$ws.Range("AF17").Value = 13.5
$ws.Range("AH17").Value = 28
$ws.Range("AI17").Value = 18.5
$ws.Range("AJ17").Value = 23
$ws.Range("H17").Value = 3.35
$ws.Range("I17").Value = 2.9
$ws.Range("N17").Value = 1.72
$ws.Range("O17").Value = 1.88
$ws.Range("T17").Value = 7.3
$ws.Range("U17").Value = 9.25
$ws.Range("V17").Value = 7.5
$ws.Range("W17").Value = 16.5
$ws.Range("Y17").Value = 19.5
$ws.Range("Z17").Value = 11
# Row 18
$ws.Range("AA18").Value = 6
$ws.Range("AB18").Value = 15
$ws.Range("AC18").Value = 51
$ws.Range("AD18").Value = 301
$ws.Range("AE18").Value = 10
$ws.Range("AF18").Value = 17
$ws.Range("AG18").Value = 13
$ws.Range("AJ18").Value = 41
$ws.Range("H18").Value = 3.2
$ws.Range("I18").Value = 3.6
$ws.Range("J18").Value = 1.07
$ws.Range("K18").Value = 8.5
$ws.Range("L18").Value = 1.33
$ws.Range("M18").Value = 3.25
$ws.Range("N18").Value = 2.1
$ws.Range("O18").Value = 1.7
$ws.Range("P18").Value = 1.44
$ws.Range("Q18").Value = 2.63
$ws.Range("R18").Value = 1.91
$ws.Range("S18").Value = 1.91
$ws.Range("T18").Value = 7
$ws.Range("U18").Value = 9.5
$ws.Range("X18").Value = 19
$ws.Range("Y18").Value = 29
$ws.Range("Z18").Value = 8.5
# Row 19
$ws.Range("AA19").Value = 8.5
$ws.Range("AD19").Value = 451
$ws.Range("AE19").Value = 15
$ws.Range("AF19").Value = 34
$ws.Range("AG19").Value = 21
$ws.Range("AH19").Value = 81
$ws.Range("AI19").Value = 51
$ws.Range("G19").Value = 1.45
$ws.Range("H19").Value = 4.5
$ws.Range("I19").Value = 6.5
$ws.Range("J19").Value = 1.04
$ws.Range("K19").Value = 13
$ws.Range("N19").Value = 1.93
$ws.Range("O19").Value = 1.93
$ws.Range("P19").Value = 1.36
$ws.Range("Q19").Value = 3
$ws.Range("R19").Value = 2.05
$ws.Range("S19").Value = 1.7
$ws.Range("T19").Value = 6
$ws.Range("U19").Value = 6.5
$ws.Range("W19").Value = 9.5
# Row 20
$ws.Range("AA20").Value = 6.5
$ws.Range("AB20").Value = 13
$ws.Range("AC20").Value = 41
$ws.Range("AD20").Value = 201
$ws.Range("AE20").Value = 9
$ws.Range("AF20").Value = 13
$ws.Range("AI20").Value = 21
$ws.Range("AJ20").Value = 29
$ws.Range("G20").Value = 2.63
$ws.Range("H20").Value = 3.4
$ws.Range("I20").Value = 2.55
$ws.Range("J20").Value = 1.05
$ws.Range("K20").Value = 11
$ws.Range("L20").Value = 1.25
$ws.Range("M20").Value = 3.75
$ws.Range("N20").Value = 1.85
$ws.Range("O20").Value = 1.95
$ws.Range("P20").Value = 1.4
$ws.Range("Q20").Value = 2.75
$ws.Range("R20").Value = 1.73
$ws.Range("S20").Value = 2
$ws.Range("T20").Value = 9.5
$ws.Range("V20").Value = 10
$ws.Range("W20").Value = 26
$ws.Range("X20").Value = 21
$ws.Range("Y20").Value = 29
$ws.Range("Z20").Value = 11
# Row 21
$ws.Range("AE21").Value = 6
$ws.Range("AF21").Value = 9.5
$ws.Range("G21").Value = 3.4
$ws.Range("I21").Value = 2.25
$ws.Range("W21").Value = 41
$ws.Range("X21").Value = 34
